$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: rename the wrap-text header cell ---
$ws.Range("A1").Value = "Initial Wrap Text "

# --- Row 7: rename "High Cell" -> "Custom Height" and make it wrap like the header cells ---
$ws.Range("A7").Value = "Custom Height"
$ws.Range("A7").WrapText = $true

# --- New row 9: "Multiple" example row ---
$ws.Range("A9").Value = "Multiple"
$ws.Range("A9").Font.Bold = $true
$ws.Range("B9").Value = "loooooong text"
$ws.Range("C9").Value = "looooong text"
$ws.Range("D9").Value = "veeeeeeeeeeeeeery loooooooooooong text"

# --- Column A: give it a fixed width ---
$ws.Columns.Item(1).ColumnWidth = 9.5

# --- Selection moves to A10 ---
[void]$ws.Range("A10").Select()

# --- Theme: lighten the window background color ---
$wb.Theme.ThemeColorScheme.Colors(2).RGB = 16777215
